# Update workbook with new TPM-derived values (rows 2-21, columns A-T)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = New-Object 'object[,]' 20,20
$data[0,0] = "ECs"
$data[0,1] = "Icam5"
$data[0,2] = "Itgal"
$data[0,3] = "ECs"
$data[0,4] = 3.0
$data[0,5] = 1.0
$data[0,6] = 1.891505666666667
$data[0,7] = 5.674517
$data[0,8] = 0.3964089626089911
$data[0,9] = 0.3964089626089912
$data[0,10] = 3.0
$data[0,11] = 1.0
$data[0,12] = 0.3547066666666667
$data[0,13] = 1.06412
$data[0,14] = 0.006934574868045491
$data[0,15] = 0.00693457486804549
$data[0,16] = 0.6709296700044444
$data[0,17] = 6.03836703004
$data[0,18] = 0.002748927629576295
$data[0,19] = 0.002748927629576295
$data[1,0] = "ECs"
$data[1,1] = "Icam5"
$data[1,2] = "Itgal"
$data[1,3] = "FAPs"
$data[1,4] = 3.0
$data[1,5] = 1.0
$data[1,6] = 1.891505666666667
$data[1,7] = 5.674517
$data[1,8] = 0.3964089626089911
$data[1,9] = 0.3964089626089912
$data[1,10] = 3.0
$data[1,11] = 1.0
$data[1,12] = 0.519977
$data[1,13] = 1.559931
$data[1,14] = 0.01016563762403213
$data[1,15] = 0.01016563762403213
$data[1,16] = 0.9835394420363333
$data[1,17] = 8.851854978327
$data[1,18] = 0.004029749864801506
$data[1,19] = 0.004029749864801506
$data[2,0] = "ECs"
$data[2,1] = "Icam5"
$data[2,2] = "Itgal"
$data[2,3] = "Inflammatory-Mac"
$data[2,4] = 3.0
$data[2,5] = 1.0
$data[2,6] = 1.891505666666667
$data[2,7] = 5.674517
$data[2,8] = 0.3964089626089911
$data[2,9] = 0.3964089626089912
$data[2,10] = 3.0
$data[2,11] = 1.0
$data[2,12] = 26.50170333333334
$data[2,13] = 79.50511
$data[2,14] = 0.5181127482682332
$data[2,15] = 0.5181127482682331
$data[2,16] = 50.12812203131889
$data[2,17] = 451.15309828187
$data[2,18] = 0.2053845370555037
$data[2,19] = 0.2053845370555037
$data[3,0] = "ECs"
$data[3,1] = "Icam5"
$data[3,2] = "Itgal"
$data[3,3] = "MuSCs"
$data[3,4] = 3.0
$data[3,5] = 1.0
$data[3,6] = 1.891505666666667
$data[3,7] = 5.674517
$data[3,8] = 0.3964089626089911
$data[3,9] = 0.3964089626089912
$data[3,10] = 3.0
$data[3,11] = 1.0
$data[3,12] = 0.09443866666666667
$data[3,13] = 0.283316
$data[3,14] = 0.001846291784117559
$data[3,15] = 0.001846291784117558
$data[3,16] = 0.1786312731524445
$data[3,17] = 1.607681458372
$data[3,18] = 0.0007318866108155448
$data[3,19] = 0.0007318866108155448
$data[4,0] = "ECs"
$data[4,1] = "Icam5"
$data[4,2] = "Itgal"
$data[4,3] = "Resolving-Mac"
$data[4,4] = 3.0
$data[4,5] = 1.0
$data[4,6] = 1.891505666666667
$data[4,7] = 5.674517
$data[4,8] = 0.3964089626089911
$data[4,9] = 0.3964089626089912
$data[4,10] = 3.0
$data[4,11] = 1.0
$data[4,12] = 23.67963033333333
$data[4,13] = 71.03889099999999
$data[4,14] = 0.4629407474555717
$data[4,15] = 0.4629407474555717
$data[4,16] = 44.79015496007188
$data[4,17] = 403.1113946406469
$data[4,18] = 0.1835138614482941
$data[4,19] = 0.1835138614482941
$data[5,0] = "FAPs"
$data[5,1] = "Icam5"
$data[5,2] = "Itgal"
$data[5,3] = "ECs"
$data[5,4] = 3.0
$data[5,5] = 1.0
$data[5,6] = 1.358655
$data[5,7] = 4.075965
$data[5,8] = 0.2847377243350503
$data[5,9] = 0.2847377243350503
$data[5,10] = 3.0
$data[5,11] = 1.0
$data[5,12] = 0.3547066666666667
$data[5,13] = 1.06412
$data[5,14] = 0.006934574868045491
$data[5,15] = 0.00693457486804549
$data[5,16] = 0.4819239862
$data[5,17] = 4.3373158758
$data[5,18] = 0.001974535067158305
$data[5,19] = 0.001974535067158305
$data[6,0] = "FAPs"
$data[6,1] = "Icam5"
$data[6,2] = "Itgal"
$data[6,3] = "FAPs"
$data[6,4] = 3.0
$data[6,5] = 1.0
$data[6,6] = 1.358655
$data[6,7] = 4.075965
$data[6,8] = 0.2847377243350503
$data[6,9] = 0.2847377243350503
$data[6,10] = 3.0
$data[6,11] = 1.0
$data[6,12] = 0.519977
$data[6,13] = 1.559931
$data[6,14] = 0.01016563762403213
$data[6,15] = 0.01016563762403213
$data[6,16] = 0.706469350935
$data[6,17] = 6.358224158415
$data[6,18] = 0.002894540523481677
$data[6,19] = 0.002894540523481677
$data[7,0] = "FAPs"
$data[7,1] = "Icam5"
$data[7,2] = "Itgal"
$data[7,3] = "Inflammatory-Mac"
$data[7,4] = 3.0
$data[7,5] = 1.0
$data[7,6] = 1.358655
$data[7,7] = 4.075965
$data[7,8] = 0.2847377243350503
$data[7,9] = 0.2847377243350503
$data[7,10] = 3.0
$data[7,11] = 1.0
$data[7,12] = 26.50170333333334
$data[7,13] = 79.50511
$data[7,14] = 0.5181127482682332
$data[7,15] = 0.5181127482682331
$data[7,16] = 36.00667174235
$data[7,17] = 324.06004568115
$data[7,18] = 0.1475262448908755
$data[7,19] = 0.1475262448908755
$data[8,0] = "FAPs"
$data[8,1] = "Icam5"
$data[8,2] = "Itgal"
$data[8,3] = "MuSCs"
$data[8,4] = 3.0
$data[8,5] = 1.0
$data[8,6] = 1.358655
$data[8,7] = 4.075965
$data[8,8] = 0.2847377243350503
$data[8,9] = 0.2847377243350503
$data[8,10] = 3.0
$data[8,11] = 1.0
$data[8,12] = 0.09443866666666667
$data[8,13] = 0.283316
$data[8,14] = 0.001846291784117559
$data[8,15] = 0.001846291784117558
$data[8,16] = 0.12830956666
$data[8,17] = 1.15478609994
$data[8,18] = 0.0005257089210681336
$data[8,19] = 0.0005257089210681336
$data[9,0] = "FAPs"
$data[9,1] = "Icam5"
$data[9,2] = "Itgal"
$data[9,3] = "Resolving-Mac"
$data[9,4] = 3.0
$data[9,5] = 1.0
$data[9,6] = 1.358655
$data[9,7] = 4.075965
$data[9,8] = 0.2847377243350503
$data[9,9] = 0.2847377243350503
$data[9,10] = 3.0
$data[9,11] = 1.0
$data[9,12] = 23.67963033333333
$data[9,13] = 71.03889099999999
$data[9,14] = 0.4629407474555717
$data[9,15] = 0.4629407474555717
$data[9,16] = 32.172448150535
$data[9,17] = 289.552033354815
$data[9,18] = 0.1318166949324667
$data[9,19] = 0.1318166949324667
$data[10,0] = "Inflammatory-Mac"
$data[10,1] = "Icam5"
$data[10,2] = "Itgal"
$data[10,3] = "ECs"
$data[10,4] = 3.0
$data[10,5] = 1.0
$data[10,6] = 0.813809
$data[10,7] = 2.441427
$data[10,8] = 0.1705525852430403
$data[10,9] = 0.1705525852430403
$data[10,10] = 3.0
$data[10,11] = 1.0
$data[10,12] = 0.3547066666666667
$data[10,13] = 1.06412
$data[10,14] = 0.006934574868045491
$data[10,15] = 0.00693457486804549
$data[10,16] = 0.2886634776933333
$data[10,17] = 2.59797129924
$data[10,18] = 0.001182709671306574
$data[10,19] = 0.001182709671306574
$data[11,0] = "Inflammatory-Mac"
$data[11,1] = "Icam5"
$data[11,2] = "Itgal"
$data[11,3] = "FAPs"
$data[11,4] = 3.0
$data[11,5] = 1.0
$data[11,6] = 0.813809
$data[11,7] = 2.441427
$data[11,8] = 0.1705525852430403
$data[11,9] = 0.1705525852430403
$data[11,10] = 3.0
$data[11,11] = 1.0
$data[11,12] = 0.519977
$data[11,13] = 1.559931
$data[11,14] = 0.01016563762403213
$data[11,15] = 0.01016563762403213
$data[11,16] = 0.423161962393
$data[11,17] = 3.808457661537
$data[11,18] = 0.001733775777422598
$data[11,19] = 0.001733775777422598
$data[12,0] = "Inflammatory-Mac"
$data[12,1] = "Icam5"
$data[12,2] = "Itgal"
$data[12,3] = "Inflammatory-Mac"
$data[12,4] = 3.0
$data[12,5] = 1.0
$data[12,6] = 0.813809
$data[12,7] = 2.441427
$data[12,8] = 0.1705525852430403
$data[12,9] = 0.1705525852430403
$data[12,10] = 3.0
$data[12,11] = 1.0
$data[12,12] = 26.50170333333334
$data[12,13] = 79.50511
$data[12,14] = 0.5181127482682332
$data[12,15] = 0.5181127482682331
$data[12,16] = 21.56732468799667
$data[12,17] = 194.10592219197
$data[12,18] = 0.08836546866452373
$data[12,19] = 0.08836546866452372
$data[13,0] = "Inflammatory-Mac"
$data[13,1] = "Icam5"
$data[13,2] = "Itgal"
$data[13,3] = "MuSCs"
$data[13,4] = 3.0
$data[13,5] = 1.0
$data[13,6] = 0.813809
$data[13,7] = 2.441427
$data[13,8] = 0.1705525852430403
$data[13,9] = 0.1705525852430403
$data[13,10] = 3.0
$data[13,11] = 1.0
$data[13,12] = 0.09443866666666667
$data[13,13] = 0.283316
$data[13,14] = 0.001846291784117559
$data[13,15] = 0.001846291784117558
$data[13,16] = 0.07685503688133334
$data[13,17] = 0.6916953319320001
$data[13,18] = 0.0003148898368942349
$data[13,19] = 0.0003148898368942349
$data[14,0] = "Inflammatory-Mac"
$data[14,1] = "Icam5"
$data[14,2] = "Itgal"
$data[14,3] = "Resolving-Mac"
$data[14,4] = 3.0
$data[14,5] = 1.0
$data[14,6] = 0.813809
$data[14,7] = 2.441427
$data[14,8] = 0.1705525852430403
$data[14,9] = 0.1705525852430403
$data[14,10] = 3.0
$data[14,11] = 1.0
$data[14,12] = 23.67963033333333
$data[14,13] = 71.03889099999999
$data[14,14] = 0.4629407474555717
$data[14,15] = 0.4629407474555717
$data[14,16] = 19.27069628193967
$data[14,17] = 173.436266537457
$data[14,18] = 0.07895574129289319
$data[14,19] = 0.0789557412928932
$data[15,0] = "Resolving-Mac"
$data[15,1] = "Icam5"
$data[15,2] = "Itgal"
$data[15,3] = "ECs"
$data[15,4] = 3.0
$data[15,5] = 1.0
$data[15,6] = 0.7076319999999999
$data[15,7] = 2.122896
$data[15,8] = 0.1483007278129181
$data[15,9] = 0.1483007278129182
$data[15,10] = 3.0
$data[15,11] = 1.0
$data[15,12] = 0.3547066666666667
$data[15,13] = 1.06412
$data[15,14] = 0.006934574868045491
$data[15,15] = 0.00693457486804549
$data[15,16] = 0.2510017879466667
$data[15,17] = 2.25901609152
$data[15,18] = 0.001028402500004317
$data[15,19] = 0.001028402500004317
$data[16,0] = "Resolving-Mac"
$data[16,1] = "Icam5"
$data[16,2] = "Itgal"
$data[16,3] = "FAPs"
$data[16,4] = 3.0
$data[16,5] = 1.0
$data[16,6] = 0.7076319999999999
$data[16,7] = 2.122896
$data[16,8] = 0.1483007278129181
$data[16,9] = 0.1483007278129182
$data[16,10] = 3.0
$data[16,11] = 1.0
$data[16,12] = 0.519977
$data[16,13] = 1.559931
$data[16,14] = 0.01016563762403213
$data[16,15] = 0.01016563762403213
$data[16,16] = 0.367952364464
$data[16,17] = 3.311571280176
$data[16,18] = 0.001507571458326349
$data[16,19] = 0.001507571458326349
$data[17,0] = "Resolving-Mac"
$data[17,1] = "Icam5"
$data[17,2] = "Itgal"
$data[17,3] = "Inflammatory-Mac"
$data[17,4] = 3.0
$data[17,5] = 1.0
$data[17,6] = 0.7076319999999999
$data[17,7] = 2.122896
$data[17,8] = 0.1483007278129181
$data[17,9] = 0.1483007278129182
$data[17,10] = 3.0
$data[17,11] = 1.0
$data[17,12] = 26.50170333333334
$data[17,13] = 79.50511
$data[17,14] = 0.5181127482682332
$data[17,15] = 0.5181127482682331
$data[17,16] = 18.75345333317333
$data[17,17] = 168.78107999856
$data[17,18] = 0.07683649765733022
$data[17,19] = 0.07683649765733022
$data[18,0] = "Resolving-Mac"
$data[18,1] = "Icam5"
$data[18,2] = "Itgal"
$data[18,3] = "MuSCs"
$data[18,4] = 3.0
$data[18,5] = 1.0
$data[18,6] = 0.7076319999999999
$data[18,7] = 2.122896
$data[18,8] = 0.1483007278129181
$data[18,9] = 0.1483007278129182
$data[18,10] = 3.0
$data[18,11] = 1.0
$data[18,12] = 0.09443866666666667
$data[18,13] = 0.283316
$data[18,14] = 0.001846291784117559
$data[18,15] = 0.001846291784117558
$data[18,16] = 0.06682782257066666
$data[18,17] = 0.601450403136
$data[18,18] = 0.0002738064153396451
$data[18,19] = 0.0002738064153396451
$data[19,0] = "Resolving-Mac"
$data[19,1] = "Icam5"
$data[19,2] = "Itgal"
$data[19,3] = "Resolving-Mac"
$data[19,4] = 3.0
$data[19,5] = 1.0
$data[19,6] = 0.7076319999999999
$data[19,7] = 2.122896
$data[19,8] = 0.1483007278129181
$data[19,9] = 0.1483007278129182
$data[19,10] = 3.0
$data[19,11] = 1.0
$data[19,12] = 23.67963033333333
$data[19,13] = 71.03889099999999
$data[19,14] = 0.4629407474555717
$data[19,15] = 0.4629407474555717
$data[19,16] = 16.75646417203733
$data[19,17] = 150.808177548336
$data[19,18] = 0.06865444978191762
$data[19,19] = 0.06865444978191762

$ws.Range("A2:T21").Value = $data
